# Fruta / hortaliza, semanal
# Insert a new weekly record as row 384 (pushing the existing rows 384-417
# down to 385-418) on the "Feria Lagunitas de Puerto Montt - Mango" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 384; everything currently at/after row 384 shifts down one.
$ws.Rows.Item(384).Insert()

# Fill in the new row 384 with the latest weekly price record.
$ws.Cells.Item(384, 1).Value = 4
$ws.Cells.Item(384, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(384, 3).Value = "Los Lagos"
$ws.Cells.Item(384, 4).Value = 45132
$ws.Cells.Item(384, 5).Value = 10
$ws.Cells.Item(384, 6).Value = "Fruta"
$ws.Cells.Item(384, 7).Value = 100108
$ws.Cells.Item(384, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(384, 9).Value = 100108002
$ws.Cells.Item(384, 10).Value = "Mango"
$ws.Cells.Item(384, 11).Value = "Sin especificar"
$ws.Cells.Item(384, 12).Value = "Primera"
$ws.Cells.Item(384, 13).Value = 100
$ws.Cells.Item(384, 14).Value = 10000
$ws.Cells.Item(384, 15).Value = 10000
$ws.Cells.Item(384, 16).Value = 10000
$ws.Cells.Item(384, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(384, 18).Value = "Perú"
$ws.Cells.Item(384, 19).Value = 2500
$ws.Cells.Item(384, 20).Value = 4
